# Studija Izvedivosti - Online Videoteka
# Update the "Analiza troskova" (Cost analysis) sheet: Ljudski resursi
# (Human resources) table gets new role names / rates, two rows are
# cleared, and a stray formatted cell (K19) is fully cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Ljudski resursi table (rows 5-11) ---

# Row 5: Projekt menadzer -> Projektni koordinator, rate 40 -> 26,
# total becomes a hard-coded formula (=6500) instead of =B5*C5
$ws.Range("A5").Value = "Projektni koordinator"
$ws.Range("C5").Value = 26
$ws.Range("D5").Formula = "=6500"

# Row 6: Frontend Developer -> Upravitelj projekta, rate 30 -> 40,
# total becomes a hard-coded value (10000) instead of the shared formula
$ws.Range("A6").Value = "Upravitelj projekta"
$ws.Range("C6").Value = 40
$ws.Range("D6").Value = 10000

# Row 7: Backend Developer -> Analiticar sustava, rate 30 -> 32,
# total becomes a hard-coded value (8000)
$ws.Range("A7").Value = "Analitičar sustava"
$ws.Range("C7").Value = 32
$ws.Range("D7").Value = 8000

# Row 8: Graficki dizajner -> Programer, total becomes a hard-coded value (5000)
$ws.Range("A8").Value = "Programer"
$ws.Range("D8").Value = 5000

# Row 9: Testni inzenjer -> Administrator baze podataka, quantity 150 -> 160,
# rate becomes the text "26,25", total becomes a hard-coded value (4200)
$ws.Range("A9").Value = "Administrator baze podataka"
$ws.Range("B9").Value = 160
$ws.Range("C9").Value = "26,25"
$ws.Range("D9").Value = 4200

# Rows 10-11 (Administrator baze podataka / Pisac dokumentacije) are no
# longer needed as separate rows - clear their contents but keep formatting
$ws.Range("A10:D10").ClearContents()
$ws.Range("A11:D11").ClearContents()

# K19 was a stray formatted-but-empty cell; fully clear it (contents +
# formatting) so it disappears from the sheet's used range entirely
$ws.Range("K19").Clear()

# Restore the view: sheet 2 active, selection on F8
$ws.Activate()
$ws.Range("F8").Select()
